$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "00000449"
$ws.Range("B9").Value = "565656"

$ws.Range("J8").Select()
